$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CasesByDate")

# --- Update existing "new cases" (column C) figures that were revised -----
# These are corrections to previously-entered daily case counts; the
# cumulative total (B) and the rolling 7-day average (D) are driven by
# formulas already in the sheet, so Excel will recompute them automatically.
$ws.Range("C266").Value2 = 1079
$ws.Range("C270").Value2 = 1226
$ws.Range("C275").Value2 = 1443
$ws.Range("C276").Value2 = 1386
$ws.Range("C283").Value2 = 2417
$ws.Range("C288").Value2 = 2824
$ws.Range("C290").Value2 = 2993
$ws.Range("C291").Value2 = 2609
$ws.Range("C295").Value2 = 3137
$ws.Range("C297").Value2 = 3002
$ws.Range("C298").Value2 = 2863
$ws.Range("C299").Value2 = 1772
$ws.Range("C300").Value2 = 1188
$ws.Range("C301").Value2 = 3585
$ws.Range("C302").Value2 = 3780
$ws.Range("C303").Value2 = 2931
$ws.Range("C304").Value2 = 430
$ws.Range("C305").Value2 = 3220
$ws.Range("C306").Value2 = 2706
$ws.Range("C307").Value2 = 1430
$ws.Range("C308").Value2 = 3804
$ws.Range("C309").Value2 = 3290

# --- Append the new day of data: Dec 2, 2020 (row 310) --------------------
# Copy the date formatting from the row above so we don't introduce a new
# number-format style, then fill in the values/formulas for the new row.
$ws.Range("A309").Copy()
$ws.Range("A310").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A310").Value2 = 44167
$ws.Range("B310").Formula = "=C310+B309"
$ws.Range("C310").Value2 = 407
$ws.Range("D310").Formula = "=AVERAGE(C304:C310)"

# --- Extend the frozen-pane selection / dimension to include the new row --
$ws.Range("C2:C310").Select()
